# Weekly cryptos-list data refresh (GitHub Actions bot commit).
# Updates the Price (D) and Volume(1h) (E) columns, plus three rows
# (36/37 and 42/44) whose coin entries swapped rank/position.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of new Price values are numeric-looking strings (e.g. "117.20",
# "1.00") that must be preserved verbatim as text, matching the source feed.
# Writing them straight to .Value lets Excel auto-convert them to numbers and
# silently drop the trailing zero (117.20 -> 117.2, 1.00 -> 1), so for those
# cells we briefly force Text format, assign the value, then restore the
# cell's original style so no formatting is actually changed.
$textValueCells = "D5","D10","D11","D12","D17","D19","D20","D23","D24","D26","D27","D28","D31","D32","D35","D36","D37","D38","D39","D41","D42","D43","D44","D49","D50","D51"
$savedStyles = @{}
foreach ($addr in $textValueCells) {
    $savedStyles[$addr] = $ws.Range($addr).Style
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "51.817.24"
$ws.Range("E2").Value = "  +4.37%  "
$ws.Range("D3").Value = "2.781.14"
$ws.Range("E3").Value = "  +5.45%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "117.20"
$ws.Range("E5").Value = "  +4.26%  "
$ws.Range("E6").Value = "  +2.87%  "
$ws.Range("E7").Value = "  +2.66%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +5.86%  "
$ws.Range("D10").Value = "42.34"
$ws.Range("E10").Value = "  +6.60%  "
$ws.Range("D11").Value = "0.0854"
$ws.Range("E11").Value = "  +5.47%  "
$ws.Range("D12").Value = "20.28"
$ws.Range("E12").Value = "  +2.21%  "
$ws.Range("E13").Value = "  +2.10%  "
$ws.Range("E14").Value = "  +4.13%  "
$ws.Range("D15").Value = "3.211.28"
$ws.Range("E15").Value = "  +5.34%  "
$ws.Range("D16").Value = "2.762.47"
$ws.Range("E16").Value = "  +5.04%  "
$ws.Range("D17").Value = "0.890"
$ws.Range("E17").Value = "  +4.63%  "
$ws.Range("D18").Value = "51.842.88"
$ws.Range("E18").Value = "  +4.70%  "
$ws.Range("D19").Value = "3.28"
$ws.Range("E19").Value = "  +10.99%  "
$ws.Range("D20").Value = "13.57"
$ws.Range("E20").Value = "  +5.18%  "
$ws.Range("E21").Value = "  +2.83%  "
$ws.Range("E22").Value = "  +3.03%  "
$ws.Range("D23").Value = "280.80"
$ws.Range("E23").Value = "  +3.85%  "
$ws.Range("D24").Value = "69.92"
$ws.Range("E24").Value = "  +1.38%  "
$ws.Range("E25").Value = "  +6.52%  "
$ws.Range("D26").Value = "26.88"
$ws.Range("E26").Value = "  +2.53%  "
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").Value = "10.23"
$ws.Range("E28").Value = "  -0.61%  "
$ws.Range("E29").Value = "  +0.92%  "
$ws.Range("E30").Value = "  +3.52%  "
$ws.Range("D31").Value = "35.34"
$ws.Range("E31").Value = "  +1.78%  "
$ws.Range("D32").Value = "50.26"
$ws.Range("E32").Value = "  +1.61%  "
$ws.Range("E33").Value = "  +1.89%  "
$ws.Range("E34").Value = "  +1.14%  "
$ws.Range("D35").Value = "19.22"
$ws.Range("E35").Value = "  +0.95%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "2.11"
$ws.Range("E36").Value = "  +3.41%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").Value = "5.00"
$ws.Range("E38").Value = "  +2.15%  "
$ws.Range("D39").Value = "3.28"
$ws.Range("E39").Value = "  +5.71%  "
$ws.Range("E40").Value = "  +11.07%  "
$ws.Range("D41").Value = "127.58"
$ws.Range("E41").Value = "  -0.44%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "2.34"
$ws.Range("E42").Value = "  +8.65%  "
$ws.Range("D43").Value = "23.39"
$ws.Range("E43").Value = "  +5.56%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "2.54"
$ws.Range("E44").Value = "  +19.69%  "
$ws.Range("E45").Value = "  +2.98%  "
$ws.Range("D46").Value = "2.091.19"
$ws.Range("E46").Value = "  +1.72%  "
$ws.Range("E47").Value = "  +3.87%  "
$ws.Range("D49").Value = "5.53"
$ws.Range("E49").Value = "  +6.16%  "
$ws.Range("D50").Value = "60.87"
$ws.Range("E50").Value = "  +2.89%  "
$ws.Range("D51").Value = "8.85"
$ws.Range("E51").Value = "  -0.58%  "

# Restore original cell styles now that the text values are safely stored.
foreach ($addr in $textValueCells) {
    $ws.Range($addr).Style = $savedStyles[$addr]
}
